# Rebuild the "concepto / monto" post-processing sheet with the new,
# bank-statement-oriented concept list (rows 2:41) instead of the old
# astor-oriented summary categories (rows 2:14).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'Impuesto Débitos y Créditos/DB'
$ws.Cells.Item(2, 2).Value = 2518822.23

$ws.Cells.Item(3, 1).Value = 'Devolución Imp. Débitos'
$ws.Cells.Item(3, 2).Value = 4729.06

$ws.Cells.Item(4, 1).Value = 'COMIS.TRANSFERENCIAS'
$ws.Cells.Item(4, 2).Value = 2530

$ws.Cells.Item(5, 1).Value = 'IVA'
$ws.Cells.Item(5, 2).Value = 57780.07999999997

$ws.Cells.Item(6, 1).Value = 'Comisión Consulta Cámara'
$ws.Cells.Item(6, 2).Value = 700

$ws.Cells.Item(7, 1).Value = 'Comisión Mantenimiento Cuenta'
$ws.Cells.Item(7, 2).Value = 4706

$ws.Cells.Item(8, 1).Value = 'Percepción I.V.A. RG. 3337'
$ws.Cells.Item(8, 2).Value = 8103.85

$ws.Cells.Item(9, 1).Value = 'Cobro Percepción IIBB'
$ws.Cells.Item(9, 2).Value = 536.72

$ws.Cells.Item(10, 1).Value = 'Débito Comisión Pago a Prov.'
$ws.Cells.Item(10, 2).Value = 784.08

$ws.Cells.Item(11, 1).Value = 'Trf. Pago.Prov-Terceros O/Bcos'
$ws.Cells.Item(11, 2).Value = 12187554.96

$ws.Cells.Item(12, 1).Value = 'Impuesto Débitos y Créditos/CR'
$ws.Cells.Item(12, 2).Value = 926362.7000000002

$ws.Cells.Item(13, 1).Value = 'IIBB- Acreditaciones Bancarias'
$ws.Cells.Item(13, 2).Value = 4150122.97

$ws.Cells.Item(14, 1).Value = 'Debito Transf. HomeBanking'
$ws.Cells.Item(14, 2).Value = 20750000

$ws.Cells.Item(15, 1).Value = 'CRED BCA ELECTR INTERBANC EXEN'
$ws.Cells.Item(15, 2).Value = 204917000

$ws.Cells.Item(16, 1).Value = 'Transferencia por CBU'
$ws.Cells.Item(16, 2).Value = 100501494.83

$ws.Cells.Item(17, 1).Value = 'CRED BCA ELECTRONICA INTERBANC'
$ws.Cells.Item(17, 2).Value = 103929174.48

$ws.Cells.Item(18, 1).Value = 'Rech. Cheques Falla Técnica'
$ws.Cells.Item(18, 2).Value = 785580.52

$ws.Cells.Item(19, 1).Value = 'Trf. Masivas Pago Proveedores'
$ws.Cells.Item(19, 2).Value = 119205508.5099999

$ws.Cells.Item(20, 1).Value = 'Pago Cheque de Cámara Recibida'
$ws.Cells.Item(20, 2).Value = 195885949.84

$ws.Cells.Item(21, 1).Value = 'SIPAP - Pago Cheque de Cámara'
$ws.Cells.Item(21, 2).Value = 32240355.95

$ws.Cells.Item(22, 1).Value = 'SIPAP - Pago Cámara SPV 24 hs.'
$ws.Cells.Item(22, 2).Value = 1592003.24

$ws.Cells.Item(23, 1).Value = 'Crédito Transf. HomeBanking'
$ws.Cells.Item(23, 2).Value = 14600000

$ws.Cells.Item(24, 1).Value = 'Descto. Docum.- Acreditación'
$ws.Cells.Item(24, 2).Value = 76168604.92000002

$ws.Cells.Item(25, 1).Value = 'Pago Cámara SPV 24 hs.'
$ws.Cells.Item(25, 2).Value = 2115162.71

$ws.Cells.Item(26, 1).Value = 'Crédito por Transferencia'
$ws.Cells.Item(26, 2).Value = 89556032.45

$ws.Cells.Item(27, 1).Value = 'Préstamos Inversion Productiva'
$ws.Cells.Item(27, 2).Value = 4455000

$ws.Cells.Item(28, 1).Value = 'DEB BCA ELECTRONICA INTERBANC'
$ws.Cells.Item(28, 2).Value = 666200.8

$ws.Cells.Item(29, 1).Value = 'Credito DEBIN'
$ws.Cells.Item(29, 2).Value = 7547000

$ws.Cells.Item(30, 1).Value = 'Pago de Servicios'
$ws.Cells.Item(30, 2).Value = 14489512.88

$ws.Cells.Item(31, 1).Value = 'Comisión Riesgo Contigente'
$ws.Cells.Item(31, 2).Value = 516

$ws.Cells.Item(32, 1).Value = 'Acreditación Cheque Dep.48 Hs.'
$ws.Cells.Item(32, 2).Value = 86020.5

$ws.Cells.Item(33, 1).Value = 'Intereses de Sobregiro'
$ws.Cells.Item(33, 2).Value = 530844.48

$ws.Cells.Item(34, 1).Value = 'IMPUESTO A LOS SELLOS'
$ws.Cells.Item(34, 2).Value = 10987.52

$ws.Cells.Item(35, 1).Value = 'Débito Automático de Servicio'
$ws.Cells.Item(35, 2).Value = 1414130.89

$ws.Cells.Item(36, 1).Value = 'Pago Automático de Préstamo'
$ws.Cells.Item(36, 2).Value = 505777.9

$ws.Cells.Item(37, 1).Value = 'DEB BCA ELECTRONICA INTRABANC'
$ws.Cells.Item(37, 2).Value = 4566105.17

$ws.Cells.Item(38, 1).Value = 'Deb. Pago de Sueldo'
$ws.Cells.Item(38, 2).Value = 8975093

$ws.Cells.Item(39, 1).Value = 'COMISIONES DATANET'
$ws.Cells.Item(39, 2).Value = 94.2

$ws.Cells.Item(40, 1).Value = 'DEB BCA ELECTR INTERBANC EXEN'
$ws.Cells.Item(40, 2).Value = 1000000

$ws.Cells.Item(41, 1).Value = 'Contras.Ints.Sobreg.'
$ws.Cells.Item(41, 2).Value = 2349.59
